# Auto-generated Excel COM-interop script
# Applies a scheduled price/profit data refresh to the Siren_Profits workbook
# (columns H-N on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 15152651
$ws.Range("I6").Value = 47619164
$ws.Range("K6").Value = 142857492
$ws.Range("M6").Value = -142857380
$ws.Range("H9").Value = 70.125
$ws.Range("J9").Value = 55
$ws.Range("L9").Value = 55
$ws.Range("N9").Value = -393
$ws.Range("H19").Value = 377.54544
$ws.Range("I19").Value = 148
$ws.Range("J19").Value = 400.5
$ws.Range("K19").Value = 148
$ws.Range("L19").Value = 400.5
$ws.Range("M19").Value = 27
$ws.Range("N19").Value = -750.5
$ws.Range("H98").Value = 16006.3545
$ws.Range("I98").Value = 17132.705
$ws.Range("K98").Value = 17132.705
$ws.Range("M98").Value = -15634.705
$ws.Range("H101").Value = 12992747
$ws.Range("I101").Value = 20414156
$ws.Range("J101").Value = 5281.75
$ws.Range("K101").Value = 61242468
$ws.Range("L101").Value = 15845.25
$ws.Range("M101").Value = -61240846
$ws.Range("N101").Value = -19089.25
$ws.Range("H106").Value = 4578459.5
$ws.Range("I106").Value = 4943536
$ws.Range("K106").Value = 4943536
$ws.Range("M106").Value = -4942905
$ws.Range("H111").Value = 772.1429000000001
$ws.Range("I111").Value = 677.1818
$ws.Range("K111").Value = 2031.5454
$ws.Range("M111").Value = 1035.4546
$ws.Range("H116").Value = 486298.97
$ws.Range("J116").Value = 4673.0835
$ws.Range("L116").Value = 4673.0835
$ws.Range("N116").Value = -11557.0835
$ws.Range("H122").Value = 16006.3545
$ws.Range("I122").Value = 17132.705
$ws.Range("K122").Value = 51398.11500000001
$ws.Range("M122").Value = -48948.11500000001
$ws.Range("H132").Value = 3660.814
$ws.Range("J132").Value = 5276
$ws.Range("L132").Value = 15828
$ws.Range("N132").Value = -20888
$ws.Range("H138").Value = 2699.3489
$ws.Range("I138").Value = 1385.4706
$ws.Range("K138").Value = 4156.4118
$ws.Range("M138").Value = 983.5882000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 10000
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -10574
$ws.Range("H61").Value = 11023.913
$ws.Range("J61").Value = 5727
$ws.Range("L61").Value = 5727
$ws.Range("N61").Value = -6151
$ws.Range("H110").Value = 1700.6086
$ws.Range("I110").Value = 1201.2354
$ws.Range("J110").Value = 3115.5
$ws.Range("K110").Value = 1201.2354
$ws.Range("L110").Value = 3115.5
$ws.Range("M110").Value = 843.7646
$ws.Range("N110").Value = -7205.5
$ws.Range("H132").Value = 3230.6572
$ws.Range("I132").Value = 1941.4348
$ws.Range("K132").Value = 5824.3044
$ws.Range("M132").Value = -3294.3044
$ws.Range("H136").Value = 11023.913
$ws.Range("J136").Value = 5727
$ws.Range("L136").Value = 17181
$ws.Range("N136").Value = -22281
$ws.Range("H139").Value = 180291.5
$ws.Range("J139").Value = 180291.5
$ws.Range("L139").Value = 180291.5
$ws.Range("N139").Value = -190571.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6336.737
$ws.Range("I86").Value = 6987.6665
$ws.Range("K86").Value = 6987.6665
$ws.Range("M86").Value = -5864.6665
$ws.Range("H89").Value = 6336.737
$ws.Range("I89").Value = 6987.6665
$ws.Range("K89").Value = 34938.3325
$ws.Range("M89").Value = -29322.3325
$ws.Range("H107").Value = 2385.4375
$ws.Range("I107").Value = 2411.1333
$ws.Range("K107").Value = 2411.1333
$ws.Range("M107").Value = -491.1333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 41154.668
$ws.Range("J50").Value = 41154.668
$ws.Range("L50").Value = 41154.668
$ws.Range("N50").Value = -42404.668
$ws.Range("H51").Value = 42155
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 42155
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 42155
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -43627
$ws.Range("H60").Value = 12953.429
$ws.Range("I60").Value = 2666.6667
$ws.Range("J60").Value = 20668.5
$ws.Range("K60").Value = 2666.6667
$ws.Range("L60").Value = 20668.5
$ws.Range("M60").Value = -2155.6667
$ws.Range("N60").Value = -21690.5
$ws.Range("H61").Value = 42155
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 42155
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 42155
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -42851
$ws.Range("H107").Value = 9375.154
$ws.Range("J107").Value = 750
$ws.Range("L107").Value = 750
$ws.Range("N107").Value = -4590
$ws.Range("H134").Value = 6333.7827
$ws.Range("I134").Value = 7391.421
$ws.Range("K134").Value = 22174.263
$ws.Range("M134").Value = -19639.263
$ws.Range("H141").Value = 268017.2
$ws.Range("J141").Value = 286399.75
$ws.Range("L141").Value = 286399.75
$ws.Range("N141").Value = -296759.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 2340.5
$ws.Range("I98").Value = 2518.8
$ws.Range("J98").Value = 2213.1428
$ws.Range("K98").Value = 7556.400000000001
$ws.Range("L98").Value = 6639.428400000001
$ws.Range("M98").Value = -6058.400000000001
$ws.Range("N98").Value = -9635.428400000001
$ws.Range("H129").Value = 23812126
$ws.Range("J129").Value = 37040536
$ws.Range("L129").Value = 111121608
$ws.Range("N129").Value = -111131608
$ws.Range("H132").Value = 38781.5
$ws.Range("I132").Value = 683
$ws.Range("J132").Value = 68413.664
$ws.Range("K132").Value = 6147
$ws.Range("L132").Value = 615722.976
$ws.Range("M132").Value = -3617
$ws.Range("N132").Value = -620782.976

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9222.083000000001
$ws.Range("J70").Value = 8807.764999999999
$ws.Range("L70").Value = 8807.764999999999
$ws.Range("N70").Value = -9347.764999999999
$ws.Range("H73").Value = 9222.083000000001
$ws.Range("J73").Value = 8807.764999999999
$ws.Range("L73").Value = 8807.764999999999
$ws.Range("N73").Value = -10679.765
$ws.Range("H80").Value = 13624
$ws.Range("I80").Value = 23770.715
$ws.Range("J80").Value = 4745.625
$ws.Range("K80").Value = 23770.715
$ws.Range("L80").Value = 4745.625
$ws.Range("M80").Value = -22772.715
$ws.Range("N80").Value = -6741.625
$ws.Range("H83").Value = 13624
$ws.Range("I83").Value = 23770.715
$ws.Range("J83").Value = 4745.625
$ws.Range("K83").Value = 118853.575
$ws.Range("L83").Value = 23728.125
$ws.Range("M83").Value = -113861.575
$ws.Range("N83").Value = -33712.125
$ws.Range("H122").Value = 12433.45
$ws.Range("I122").Value = 7998.625
$ws.Range("J122").Value = 30172.75
$ws.Range("K122").Value = 23995.875
$ws.Range("L122").Value = 90518.25
$ws.Range("M122").Value = -21545.875
$ws.Range("N122").Value = -95418.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3035.2307
$ws.Range("I82").Value = 5391.4
$ws.Range("K82").Value = 5391.4
$ws.Range("M82").Value = -5030.4
$ws.Range("H85").Value = 3035.2307
$ws.Range("I85").Value = 5391.4
$ws.Range("K85").Value = 5391.4
$ws.Range("M85").Value = -4143.4
$ws.Range("H100").Value = 6748.391
$ws.Range("I100").Value = 7722.0835
$ws.Range("J100").Value = 5686.1816
$ws.Range("K100").Value = 7722.0835
$ws.Range("L100").Value = 5686.1816
$ws.Range("M100").Value = -7181.0835
$ws.Range("N100").Value = -6768.1816
$ws.Range("H136").Value = 5957.8945
$ws.Range("I136").Value = 4179.8
$ws.Range("K136").Value = 12539.4
$ws.Range("M136").Value = -9989.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 24296.08
$ws.Range("I122").Value = 2401.5
$ws.Range("J122").Value = 34599.41
$ws.Range("K122").Value = 7204.5
$ws.Range("L122").Value = 103798.23
$ws.Range("M122").Value = -4754.5
$ws.Range("N122").Value = -108698.23
$ws.Range("H136").Value = 423086.66
$ws.Range("I136").Value = 555008.3
$ws.Range("K136").Value = 1665024.9
$ws.Range("M136").Value = -1662474.9

